# Natmi following Dr Hou advice
#
# A third cell cluster ("FAPs") is added to the sending/target cluster
# combinations analysed in this L1cam -> Ephb2 ligand-receptor pair sheet.
# This grows the cluster-pair table from 2x2 (ECs, sCs) to 3x3 minus the
# pairs targeting ECs (ECs is not receptor-expressing for Ephb2), i.e. from
# 4 data rows to 6 data rows, and refreshes every computed statistic
# (expression / specificity values) that depends on the enlarged pool of
# clusters.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Row($r, $values) {
    $cols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T")
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range($cols[$i] + $r).Value = $values[$i]
    }
}

# Row 2: ECs -> L1cam -> Ephb2 -> FAPs
Set-Row 2 @(
    "ECs", "L1cam", "Ephb2", "FAPs",
    3, 1, 19.72083766666667, 59.162513, 0.8016210077351786, 0.8016210077351787,
    3, 1, 6.346253666666667, 19.038761, 0.9446330608455225, 0.9446330608455226,
    125.1534383518214, 1126.380945166393, 0.757237706174954, 0.7572377061749542
)

# Row 3: ECs -> L1cam -> Ephb2 -> sCs
Set-Row 3 @(
    "ECs", "L1cam", "Ephb2", "sCs",
    3, 1, 19.72083766666667, 59.162513, 0.8016210077351786, 0.8016210077351787,
    3, 1, 0.3719673333333333, 1.115902, 0.05536693915447755, 0.05536693915447755,
    7.335507397969556, 66.019566581726, 0.04438330156022462, 0.04438330156022462
)

# Row 4 (new): FAPs -> L1cam -> Ephb2 -> FAPs
Set-Row 4 @(
    "FAPs", "L1cam", "Ephb2", "FAPs",
    1, 0.3333333333333333, 0.099159, 0.297477, 0.004030657259573097, 0.004030657259573097,
    3, 1, 6.346253666666667, 19.038761, 0.9446330608455225, 0.9446330608455226,
    0.629288167333, 5.663593505997, 0.00380749210432976, 0.003807492104329761
)

# Row 5 (new): FAPs -> L1cam -> Ephb2 -> sCs
Set-Row 5 @(
    "FAPs", "L1cam", "Ephb2", "sCs",
    1, 0.3333333333333333, 0.099159, 0.297477, 0.004030657259573097, 0.004030657259573097,
    3, 1, 0.3719673333333333, 1.115902, 0.05536693915447755, 0.05536693915447755,
    0.036883908806, 0.331955179254, 0.0002231651552433369, 0.0002231651552433369
)

# Row 6 (previously row 4): sCs -> L1cam -> Ephb2 -> FAPs
Set-Row 6 @(
    "sCs", "L1cam", "Ephb2", "FAPs",
    3, 1, 4.781202, 14.343606, 0.1943483350052483, 0.1943483350052483,
    3, 1, 6.346253666666667, 19.038761, 0.9446330608455225, 0.9446330608455226,
    30.342720723574, 273.0844865121661, 0.1835878625662387, 0.1835878625662387
)

# Row 7 (previously row 5): sCs -> L1cam -> Ephb2 -> sCs
Set-Row 7 @(
    "sCs", "L1cam", "Ephb2", "sCs",
    3, 1, 4.781202, 14.343606, 0.1943483350052483, 0.1943483350052483,
    3, 1, 0.3719673333333333, 1.115902, 0.05536693915447755, 0.05536693915447755,
    1.778450958068, 16.006058622612, 0.0107604724390096, 0.0107604724390096
)
